$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# "Step 3" sheet: the reagent in row 3 was re-identified. It used to be
# entered as a pre-made "3M methanolic ammonia solution" (with its
# pseudo-SMILES and a Fisher Scientific link); it's now annotated using
# the actual reagent, ammonia, with its real SMILES and a Sigma-Aldrich
# product link. The buy-mass/buy-CAD/rxn-qty numbers were updated to
# match the new vendor listing.
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Step 3")

$ws3.Range("A3").Value = "ammonia"
$ws3.Range("B3").Value = "N"
$ws3.Range("C3").Value = 17.03
$ws3.Range("D3").Value = "https://www.sigmaaldrich.com/catalog/product/aldrich/294993?lang=en&region=CA"
$ws3.Range("E3").Value = 828
$ws3.Range("F3").Value = 170

# Point the D3 hyperlink at the new vendor page instead of the old one.
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://www.sigmaaldrich.com/catalog/product/aldrich/294993?lang=en&region=CA")
$ws3.Range("D3").Style = "Hyperlink"

# "Step 3" becomes the selected/active tab, with A3:F3 (the edited row)
# highlighted.
$null = $ws3.Activate()
$null = $ws3.Range("A3:F3").Select()
